$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column range to Text format before assigning,
# so numeric-looking strings (e.g. "206.59", "1.00") are kept as text
# instead of being auto-converted to numbers, matching the original
# inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Price (column D) updates ---
$ws.Range("D2").Value = "87.007.06"
$ws.Range("D3").Value = "3.160.74"
$ws.Range("D5").Value = "206.59"
$ws.Range("D6").Value = "605.32"
$ws.Range("D7").Value = "0.364"
$ws.Range("D8").Value = "0.653"
$ws.Range("D9").Value = "1.00"
$ws.Range("D10").Value = "3.164.10"
$ws.Range("D11").Value = "0.533"
$ws.Range("D13").Value = "0.0000241"
$ws.Range("D14").Value = "3.754.48"
$ws.Range("D15").Value = "5.24"
$ws.Range("D16").Value = "86.767.48"
$ws.Range("D17").Value = "31.96"
$ws.Range("D18").Value = "3.201.64"
$ws.Range("D19").Value = "2.94"
$ws.Range("D20").Value = "13.32"
$ws.Range("D21").Value = "410.66"
$ws.Range("D22").Value = "8.44"
$ws.Range("D23").Value = "5.01"
$ws.Range("D24").Value = "5.11"
$ws.Range("D25").Value = "11.53"
$ws.Range("D26").Value = "3.344.77"
$ws.Range("D27").Value = "73.19"
$ws.Range("D29").Value = "0.999"
$ws.Range("D30").Value = "0.161"
$ws.Range("D32").Value = "536.85"
$ws.Range("D33").Value = "8.28"
$ws.Range("D35").Value = "1.27"
$ws.Range("D36").Value = "6.57"
$ws.Range("D38").Value = "21.62"
$ws.Range("D39").Value = "21.80"
$ws.Range("D40").Value = "1.00"
$ws.Range("D41").Value = "2.99"
$ws.Range("D43").Value = "0.373"
$ws.Range("D45").Value = "148.67"
$ws.Range("D46").Value = "171.46"
$ws.Range("D47").Value = "43.17"
$ws.Range("D50").Value = "3.93"
$ws.Range("D51").Value = "0.586"

# Restore default (General/no explicit) cell format so no stray
# style is left applied to the cells themselves.
$ws.Range("D2:D51").ClearFormats()

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("E3").Value = "  -5.02%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E6").Value = "  -6.71%  "
$ws.Range("E7").Value = "  -7.42%  "
$ws.Range("E8").Value = "  +8.39%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -4.88%  "
$ws.Range("E11").Value = "  -8.98%  "
$ws.Range("E12").Value = "  +5.42%  "
$ws.Range("E13").Value = "  -15.73%  "
$ws.Range("E14").Value = "  -4.80%  "
$ws.Range("E15").Value = "  -4.80%  "
$ws.Range("E16").Value = "  -1.84%  "
$ws.Range("E17").Value = "  -9.63%  "
$ws.Range("E18").Value = "  -4.06%  "
$ws.Range("E19").Value = "  -5.30%  "
$ws.Range("E20").Value = "  -8.78%  "
$ws.Range("E21").Value = "  -9.72%  "
$ws.Range("E22").Value = "  -12.67%  "
$ws.Range("E23").Value = "  -7.97%  "
$ws.Range("E24").Value = "  -7.31%  "
$ws.Range("E25").Value = "  -10.18%  "
$ws.Range("E26").Value = "  -4.78%  "
$ws.Range("E27").Value = "  -6.60%  "
$ws.Range("E28").Value = "  -3.49%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -24.16%  "
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("E32").Value = "  -9.16%  "
$ws.Range("E33").Value = "  -11.69%  "
$ws.Range("E34").Value = "  -12.63%  "
$ws.Range("E35").Value = "  -20.21%  "
$ws.Range("E36").Value = "  -7.87%  "
$ws.Range("E37").Value = "  -5.92%  "
$ws.Range("E38").Value = "  -6.37%  "
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("E41").Value = "  -4.55%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  -11.18%  "
$ws.Range("E44").Value = "  -11.92%  "
$ws.Range("E45").Value = "  -6.02%  "
$ws.Range("E46").Value = "  -8.64%  "
$ws.Range("E47").Value = "  -6.17%  "
$ws.Range("E48").Value = "  +9.63%  "
$ws.Range("E49").Value = "  -13.02%  "
$ws.Range("E50").Value = "  -11.09%  "
$ws.Range("E51").Value = "  -11.02%  "
